$wb = $excel.ActiveWorkbook

# --- Sheet1 ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Range("C2").Value = "Yes"
$ws1.Range("E2").Select()

# --- DATA sheet ---
$ws2 = $wb.Worksheets.Item("DATA")

# Row2: C2 gains quote-prefix style (value "chrome" unchanged)
$ws2.Range("C2").Formula = "'chrome"

# Row3: B3 "yes" -> "No"; C3 gains quote-prefix style
$ws2.Range("B3").Value = "No"
$ws2.Range("C3").Formula = "'chrome"

# Row4: B4 "no" -> "No"; C4 gains quote-prefix style
$ws2.Range("B4").Value = "No"
$ws2.Range("C4").Formula = "'chrome"

# Row5: B5 "no" -> "No"; C5 gains quote-prefix style
$ws2.Range("B5").Value = "No"
$ws2.Range("C5").Formula = "'chrome"

# Row6: B6 "no" -> "No"; C6 gains quote-prefix style; D6 "Admin123" -> "Admin" (loses quote-prefix style)
$ws2.Range("B6").Value = "No"
$ws2.Range("C6").Formula = "'chrome"
$ws2.Range("D6").Value = "Admin"

$ws2.Range("C6").Select()
